$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 124.25
$ws.Range("I12").Value = 99
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 99
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = 71
$ws.Range("N12").Value = -540

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 293.78946
$ws.Range("I19").Value = 313.875
$ws.Range("J19").Value = 279.18182
$ws.Range("K19").Value = 313.875
$ws.Range("L19").Value = 279.18182
$ws.Range("M19").Value = -138.875
$ws.Range("N19").Value = -629.18182

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 57514.25
$ws.Range("I21").Value = 57514.25
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 57514.25
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -57046.25
$ws.Range("N21").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 57514.25
$ws.Range("I23").Value = 57514.25
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 57514.25
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -57280.25
$ws.Range("N23").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 481.8421
$ws.Range("I28").Value = 374.92856
$ws.Range("J28").Value = 781.2
$ws.Range("K28").Value = 374.92856
$ws.Range("L28").Value = 781.2
$ws.Range("M28").Value = 110.07144
$ws.Range("N28").Value = -1751.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 778.9
$ws.Range("I58").Value = 312.7143
$ws.Range("J58").Value = 1866.6666
$ws.Range("K58").Value = 938.1428999999999
$ws.Range("L58").Value = 5599.9998
$ws.Range("M58").Value = -788.1428999999999
$ws.Range("N58").Value = -5899.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3150.5789
$ws.Range("I62").Value = 3278.3704
$ws.Range("K62").Value = 3278.3704
$ws.Range("M62").Value = -2654.3704

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3150.5789
$ws.Range("I65").Value = 3278.3704
$ws.Range("K65").Value = 16391.852
$ws.Range("M65").Value = -13271.852

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 7147331.5
$ws.Range("I86").Value = 16670640
$ws.Range("J86").Value = 4849.875
$ws.Range("K86").Value = 16670640
$ws.Range("L86").Value = 4849.875
$ws.Range("M86").Value = -16669517
$ws.Range("N86").Value = -7095.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 7147331.5
$ws.Range("I89").Value = 16670640
$ws.Range("J89").Value = 4849.875
$ws.Range("K89").Value = 83353200
$ws.Range("L89").Value = 24249.375
$ws.Range("M89").Value = -83347584
$ws.Range("N89").Value = -35481.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 510000
$ws.Range("I14").Value = 1000000
$ws.Range("J14").Value = 20000
$ws.Range("K14").Value = 1000000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = -999825
$ws.Range("N14").Value = -20350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 845.8293
$ws.Range("I97").Value = 816.54285
$ws.Range("J97").Value = 1016.6667
$ws.Range("K97").Value = 816.54285
$ws.Range("L97").Value = 1016.6667
$ws.Range("M97").Value = -320.54285
$ws.Range("N97").Value = -2008.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 37487.5
$ws.Range("J133").Value = 37487.5
$ws.Range("L133").Value = 37487.5
$ws.Range("N133").Value = -42547.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 28617.715
$ws.Range("J140").Value = 28617.715
$ws.Range("L140").Value = 28617.715
$ws.Range("N140").Value = -38977.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 882.75
$ws.Range("I22").Value = 882.75
$ws.Range("K22").Value = 882.75
$ws.Range("M22").Value = -709.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 22359.924
$ws.Range("J140").Value = 22359.924
$ws.Range("L140").Value = 22359.924
$ws.Range("N140").Value = -32719.924

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 73800.64999999999
$ws.Range("I132").Value = 2280.9092
$ws.Range("K132").Value = 6842.7276
$ws.Range("M132").Value = -4312.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5140.5654
$ws.Range("J5").Value = 1897
$ws.Range("L5").Value = 5691
$ws.Range("N5").Value = -5915

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 24.357143
$ws.Range("I14").Value = 24.357143
$ws.Range("K14").Value = 73.07142899999999
$ws.Range("M14").Value = 99.92857100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 26695158
$ws.Range("J33").Value = 33368898
$ws.Range("L33").Value = 200213388
$ws.Range("N33").Value = -200213954

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1994.0555
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1994.0555
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5982.166499999999
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -6150.166499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2336.3635
$ws.Range("J39").Value = 2550
$ws.Range("L39").Value = 7650
$ws.Range("N39").Value = -8238

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1968.1818
$ws.Range("J55").Value = 1968.1818
$ws.Range("L55").Value = 5904.5454
$ws.Range("N55").Value = -6258.5454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1460
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 1460
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 4380
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -7376

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 3994
$ws.Range("I103").Value = 7575
$ws.Range("J103").Value = 2203.5
$ws.Range("K103").Value = 22725
$ws.Range("L103").Value = 6610.5
$ws.Range("M103").Value = -21846
$ws.Range("N103").Value = -8368.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3024.82
$ws.Range("I113").Value = 4017.2415
$ws.Range("J113").Value = 1654.3334
$ws.Range("K113").Value = 12051.7245
$ws.Range("L113").Value = 4963.0002
$ws.Range("M113").Value = -9881.7245
$ws.Range("N113").Value = -9303.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 8942.385
$ws.Range("I122").Value = 424
$ws.Range("J122").Value = 22571.8
$ws.Range("K122").Value = 3816
$ws.Range("L122").Value = 203146.2
$ws.Range("M122").Value = -1366
$ws.Range("N122").Value = -208046.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4428.2573
$ws.Range("I131").Value = 11523.223
$ws.Range("J131").Value = 1972.3077
$ws.Range("K131").Value = 34569.669
$ws.Range("L131").Value = 5916.9231
$ws.Range("M131").Value = -29529.669
$ws.Range("N131").Value = -15996.9231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 5140.5654
$ws.Range("J135").Value = 1897
$ws.Range("L135").Value = 17073
$ws.Range("N135").Value = -22143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 47800
$ws.Range("J88").Value = 47800
$ws.Range("L88").Value = 47800
$ws.Range("N88").Value = -48702

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H91").Value = 47800
$ws.Range("J91").Value = 47800
$ws.Range("L91").Value = 47800
$ws.Range("N91").Value = -50920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1858.08
$ws.Range("I102").Value = 1801.7646
$ws.Range("J102").Value = 1977.75
$ws.Range("K102").Value = 1801.7646
$ws.Range("L102").Value = 1977.75
$ws.Range("M102").Value = -179.7646
$ws.Range("N102").Value = -5221.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 46626
$ws.Range("J108").Value = 46626
$ws.Range("L108").Value = 46626
$ws.Range("N108").Value = -54306

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 9899.5
$ws.Range("J45").Value = 9800
$ws.Range("L45").Value = 9800
$ws.Range("N45").Value = -10782
